{"js": "// Update the stack-trace text in the document body to match the new\n// line numbers / class names produced after moving from M2Doc 2.0.1 to\n// 2.0.2 (per the commit message). The entire stack trace lives inside a\n// single w:t run, so we locate and replace each changed fragment via\n// Word's search-and-replace (using enough surrounding context to target\n// the correct occurrence when a line number repeats elsewhere).\n\nconst replacements = [\n  {\n    find: \"M2DocEvaluator.caseUserDoc(M2DocEvaluator.java:1132)\",\n    replace: \"M2DocEvaluator.caseUserDoc(M2DocEvaluator.java:1252)\"\n  },\n  {\n    find: \"M2DocEvaluator.doSwitch(M2DocEvaluator.java:1096)\\n\\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseBlock(M2DocEvaluator.java:1305)\",\n    replace: \"M2DocEvaluator.doSwitch(M2DocEvaluator.java:1216)\\n\\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseBlock(M2DocEvaluator.java:1425)\"\n  },\n  {\n    find: \"M2DocEvaluator.doSwitch(M2DocEvaluator.java:1096)\\n\\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseDocumentTemplate(M2DocEvaluator.java:283)\",\n    replace: \"M2DocEvaluator.doSwitch(M2DocEvaluator.java:1216)\\n\\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseDocumentTemplate(M2DocEvaluator.java:287)\"\n  },\n  {\n    find: \"M2DocEvaluator.doSwitch(M2DocEvaluator.java:1096)\\n\\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.generate(M2DocEvaluator.java:272)\",\n    replace: \"M2DocEvaluator.doSwitch(M2DocEvaluator.java:1216)\\n\\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.generate(M2DocEvaluator.java:276)\"\n  },\n  {\n    find: \"AbstractTemplatesTestSuite.prepareoutputAndGenerate(AbstractTemplatesTestSuite.java:479)\\n\\tat org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.generation(AbstractTemplatesTestSuite.java:388)\\n\\tat sun.reflect.GeneratedMethodAccessor75.invoke(Unknown Source)\",\n    replace: \"AbstractTemplatesTestSuite.prepareoutputAndGenerate(AbstractTemplatesTestSuite.java:480)\\n\\tat org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.generation(AbstractTemplatesTestSuite.java:389)\\n\\tat sun.reflect.GeneratedMethodAccessor74.invoke(Unknown Source)\"\n  },\n  {\n    find: \"\\tat org.eclipse.jdt.internal.junit4.runner.JUnit4TestReference.run(JUnit4TestReference.java:86)\\n\" +\n      \"\\tat org.eclipse.jdt.internal.junit.runner.TestExecution.run(TestExecution.java:38)\\n\" +\n      \"\\tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:539)\\n\" +\n      \"\\tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:761)\\n\" +\n      \"\\tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.run(RemoteTestRunner.java:461)\\n\" +\n      \"\\tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.main(RemoteTestRunner.java:207)\",\n    replace: \"\\tat org.apache.maven.surefire.junit4.JUnit4Provider.execute(JUnit4Provider.java:264)\\n\" +\n      \"\\tat org.apache.maven.surefire.junit4.JUnit4Provider.executeTestSet(JUnit4Provider.java:153)\\n\" +\n      \"\\tat org.apache.maven.surefire.junit4.JUnit4Provider.invoke(JUnit4Provider.java:124)\\n\" +\n      \"\\tat sun.reflect.NativeMethodAccessorImpl.invoke0(Native Method)\\n\" +\n      \"\\tat sun.reflect.NativeMethodAccessorImpl.invoke(NativeMethodAccessorImpl.java:62)\\n\" +\n      \"\\tat sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)\\n\" +\n      \"\\tat java.lang.reflect.Method.invoke(Method.java:498)\\n\" +\n      \"\\tat org.apache.maven.surefire.util.ReflectionUtils.invokeMethodWithArray2(ReflectionUtils.java:208)\\n\" +\n      \"\\tat org.apache.maven.surefire.booter.ProviderFactory$ProviderProxy.invoke(ProviderFactory.java:156)\\n\" +\n      \"\\tat org.apache.maven.surefire.booter.ProviderFactory.invokeProvider(ProviderFactory.java:82)\\n\" +\n      \"\\tat org.eclipse.tycho.surefire.osgibooter.OsgiSurefireBooter.run(OsgiSurefireBooter.java:91)\\n\" +\n      \"\\tat org.eclipse.tycho.surefire.osgibooter.HeadlessTestApplication.run(HeadlessTestApplication.java:21)\\n\" +\n      \"\\tat sun.reflect.NativeMethodAccessorImpl.invoke0(Native Method)\\n\" +\n      \"\\tat sun.reflect.NativeMethodAccessorImpl.invoke(NativeMethodAccessorImpl.java:62)\\n\" +\n      \"\\tat sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)\\n\" +\n      \"\\tat java.lang.reflect.Method.invoke(Method.java:498)\\n\" +\n      \"\\tat org.eclipse.equinox.internal.app.EclipseAppContainer.callMethodWithException(EclipseAppContainer.java:587)\\n\" +\n      \"\\tat org.eclipse.equinox.internal.app.EclipseAppHandle.run(EclipseAppHandle.java:198)\\n\" +\n      \"\\tat org.eclipse.core.runtime.internal.adaptor.EclipseAppLauncher.runApplication(EclipseAppLauncher.java:134)\\n\" +\n      \"\\tat org.eclipse.core.runtime.internal.adaptor.EclipseAppLauncher.start(EclipseAppLauncher.java:104)\\n\" +\n      \"\\tat org.eclipse.core.runtime.adaptor.EclipseStarter.run(EclipseStarter.java:388)\\n\" +\n      \"\\tat org.eclipse.core.runtime.adaptor.EclipseStarter.run(EclipseStarter.java:243)\\n\" +\n      \"\\tat sun.reflect.NativeMethodAccessorImpl.invoke0(Native Method)\\n\" +\n      \"\\tat sun.reflect.NativeMethodAccessorImpl.invoke(NativeMethodAccessorImpl.java:62)\\n\" +\n      \"\\tat sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)\\n\" +\n      \"\\tat java.lang.reflect.Method.invoke(Method.java:498)\\n\" +\n      \"\\tat org.eclipse.equinox.launcher.Main.invokeFramework(Main.java:656)\\n\" +\n      \"\\tat org.eclipse.equinox.launcher.Main.basicRun(Main.java:592)\\n\" +\n      \"\\tat org.eclipse.equinox.launcher.Main.run(Main.java:1498)\\n\" +\n      \"\\tat org.eclipse.equinox.launcher.Main.main(Main.java:1471)\"\n  }\n];\n\nconst body = context.document.body;\n\nfor (const { find, replace } of replacements) {\n  const results = body.search(find, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + find.substring(0, 60));\n  }\n\n  results.items[0].insertText(replace, \"Replace\");\n  await context.sync();\n}\n", "ps1": "# Update the stack-trace text in the document to match the new line\n# numbers / class names produced after moving M2Doc from 2.0.1 to 2.0.2\n# (per the commit message). The whole stack trace lives in a single run\n# of text, so each changed fragment is located and replaced with\n# Find/Replace. Enough surrounding context (including embedded tabs and\n# newlines via backtick escapes) is included in each search string so\n# that the correct occurrence is matched when a line number repeats\n# elsewhere in the trace.\n#\n# NOTE: each Find/Replace below is issued as its own top-level\n# statement block (its own fresh Find object, not funneled through a\n# shared helper function) and is not the final statement of the\n# script - this interpreter's COM shim only reliably commits a\n# Find.Execute(Replace:=...) call made from inside a user-defined\n# function the *first* time that function is invoked; later calls\n# silently no-op. Keeping every call inlined at the top level sidesteps\n# that and reliably applies all six edits.\n\n$d = $word.ActiveDocument\n\n$find1 = $d.Content.Find\n$find1.ClearFormatting()\n$find1.Replacement.ClearFormatting()\n$result1 = $find1.Execute( `\n    \"M2DocEvaluator.caseUserDoc(M2DocEvaluator.java:1132)\", `\n    $true, $false, $false, $false, $false, $true, 1, $false, `\n    \"M2DocEvaluator.caseUserDoc(M2DocEvaluator.java:1252)\", 1)\nif (-not $result1) {\n    Write-Output \"WARNING: replacement 1 did not find a match\"\n}\n\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Replacement.ClearFormatting()\n$result2 = $find2.Execute( `\n    \"M2DocEvaluator.doSwitch(M2DocEvaluator.java:1096)`n`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseBlock(M2DocEvaluator.java:1305)\", `\n    $true, $false, $false, $false, $false, $true, 1, $false, `\n    \"M2DocEvaluator.doSwitch(M2DocEvaluator.java:1216)`n`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseBlock(M2DocEvaluator.java:1425)\", 1)\nif (-not $result2) {\n    Write-Output \"WARNING: replacement 2 did not find a match\"\n}\n\n$find3 = $d.Content.Find\n$find3.ClearFormatting()\n$find3.Replacement.ClearFormatting()\n$result3 = $find3.Execute( `\n    \"M2DocEvaluator.doSwitch(M2DocEvaluator.java:1096)`n`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseDocumentTemplate(M2DocEvaluator.java:283)\", `\n    $true, $false, $false, $false, $false, $true, 1, $false, `\n    \"M2DocEvaluator.doSwitch(M2DocEvaluator.java:1216)`n`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseDocumentTemplate(M2DocEvaluator.java:287)\", 1)\nif (-not $result3) {\n    Write-Output \"WARNING: replacement 3 did not find a match\"\n}\n\n$find4 = $d.Content.Find\n$find4.ClearFormatting()\n$find4.Replacement.ClearFormatting()\n$result4 = $find4.Execute( `\n    \"M2DocEvaluator.doSwitch(M2DocEvaluator.java:1096)`n`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.generate(M2DocEvaluator.java:272)\", `\n    $true, $false, $false, $false, $false, $true, 1, $false, `\n    \"M2DocEvaluator.doSwitch(M2DocEvaluator.java:1216)`n`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.generate(M2DocEvaluator.java:276)\", 1)\nif (-not $result4) {\n    Write-Output \"WARNING: replacement 4 did not find a match\"\n}\n\n$find5 = $d.Content.Find\n$find5.ClearFormatting()\n$find5.Replacement.ClearFormatting()\n$result5 = $find5.Execute( `\n    \"AbstractTemplatesTestSuite.prepareoutputAndGenerate(AbstractTemplatesTestSuite.java:479)`n`tat org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.generation(AbstractTemplatesTestSuite.java:388)`n`tat sun.reflect.GeneratedMethodAccessor75.invoke(Unknown Source)\", `\n    $true, $false, $false, $false, $false, $true, 1, $false, `\n    \"AbstractTemplatesTestSuite.prepareoutputAndGenerate(AbstractTemplatesTestSuite.java:480)`n`tat org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.generation(AbstractTemplatesTestSuite.java:389)`n`tat sun.reflect.GeneratedMethodAccessor74.invoke(Unknown Source)\", 1)\nif (-not $result5) {\n    Write-Output \"WARNING: replacement 5 did not find a match\"\n}\n\n$oldTail = \"`tat org.eclipse.jdt.internal.junit4.runner.JUnit4TestReference.run(JUnit4TestReference.java:86)`n\" + `\n    \"`tat org.eclipse.jdt.internal.junit.runner.TestExecution.run(TestExecution.java:38)`n\" + `\n    \"`tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:539)`n\" + `\n    \"`tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:761)`n\" + `\n    \"`tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.run(RemoteTestRunner.java:461)`n\" + `\n    \"`tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.main(RemoteTestRunner.java:207)\"\n\n$newTail = \"`tat org.apache.maven.surefire.junit4.JUnit4Provider.execute(JUnit4Provider.java:264)`n\" + `\n    \"`tat org.apache.maven.surefire.junit4.JUnit4Provider.executeTestSet(JUnit4Provider.java:153)`n\" + `\n    \"`tat org.apache.maven.surefire.junit4.JUnit4Provider.invoke(JUnit4Provider.java:124)`n\" + `\n    \"`tat sun.reflect.NativeMethodAccessorImpl.invoke0(Native Method)`n\" + `\n    \"`tat sun.reflect.NativeMethodAccessorImpl.invoke(NativeMethodAccessorImpl.java:62)`n\" + `\n    \"`tat sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)`n\" + `\n    \"`tat java.lang.reflect.Method.invoke(Method.java:498)`n\" + `\n    \"`tat org.apache.maven.surefire.util.ReflectionUtils.invokeMethodWithArray2(ReflectionUtils.java:208)`n\" + `\n    \"`tat org.apache.maven.surefire.booter.ProviderFactory`$ProviderProxy.invoke(ProviderFactory.java:156)`n\" + `\n    \"`tat org.apache.maven.surefire.booter.ProviderFactory.invokeProvider(ProviderFactory.java:82)`n\" + `\n    \"`tat org.eclipse.tycho.surefire.osgibooter.OsgiSurefireBooter.run(OsgiSurefireBooter.java:91)`n\" + `\n    \"`tat org.eclipse.tycho.surefire.osgibooter.HeadlessTestApplication.run(HeadlessTestApplication.java:21)`n\" + `\n    \"`tat sun.reflect.NativeMethodAccessorImpl.invoke0(Native Method)`n\" + `\n    \"`tat sun.reflect.NativeMethodAccessorImpl.invoke(NativeMethodAccessorImpl.java:62)`n\" + `\n    \"`tat sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)`n\" + `\n    \"`tat java.lang.reflect.Method.invoke(Method.java:498)`n\" + `\n    \"`tat org.eclipse.equinox.internal.app.EclipseAppContainer.callMethodWithException(EclipseAppContainer.java:587)`n\" + `\n    \"`tat org.eclipse.equinox.internal.app.EclipseAppHandle.run(EclipseAppHandle.java:198)`n\" + `\n    \"`tat org.eclipse.core.runtime.internal.adaptor.EclipseAppLauncher.runApplication(EclipseAppLauncher.java:134)`n\" + `\n    \"`tat org.eclipse.core.runtime.internal.adaptor.EclipseAppLauncher.start(EclipseAppLauncher.java:104)`n\" + `\n    \"`tat org.eclipse.core.runtime.adaptor.EclipseStarter.run(EclipseStarter.java:388)`n\" + `\n    \"`tat org.eclipse.core.runtime.adaptor.EclipseStarter.run(EclipseStarter.java:243)`n\" + `\n    \"`tat sun.reflect.NativeMethodAccessorImpl.invoke0(Native Method)`n\" + `\n    \"`tat sun.reflect.NativeMethodAccessorImpl.invoke(NativeMethodAccessorImpl.java:62)`n\" + `\n    \"`tat sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)`n\" + `\n    \"`tat java.lang.reflect.Method.invoke(Method.java:498)`n\" + `\n    \"`tat org.eclipse.equinox.launcher.Main.invokeFramework(Main.java:656)`n\" + `\n    \"`tat org.eclipse.equinox.launcher.Main.basicRun(Main.java:592)`n\" + `\n    \"`tat org.eclipse.equinox.launcher.Main.run(Main.java:1498)`n\" + `\n    \"`tat org.eclipse.equinox.launcher.Main.main(Main.java:1471)\"\n\n$find6 = $d.Content.Find\n$find6.ClearFormatting()\n$find6.Replacement.ClearFormatting()\n$result6 = $find6.Execute($oldTail, $true, $false, $false, $false, $false, $true, 1, $false, $newTail, 1)\nif (-not $result6) {\n    Write-Output \"WARNING: replacement 6 did not find a match\"\n}\n\nWrite-Output \"stack trace updated\"\n"}
